$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the new "SeaBattle" worksheet as the last tab (sheetId 4, rId4),
#    which also naturally bumps the workbook's activeTab to the new sheet.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws.Name = "SeaBattle"

# ---------------------------------------------------------------------------
# 2. Header row (row 1) -- written left-to-right so any brand-new shared
#    strings are interned in the same order the original commit used.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Skull total"
$ws.Range("B1").Value = "Monkey"
$ws.Range("C1").Value = "Parrot"
$ws.Range("D1").Value = "Sword"
$ws.Range("E1").Value = "Coin"
$ws.Range("F1").Value = "Diamond"
$ws.Range("G1").Value = "expect score"
$ws.Range("H1").Value = "sword require"
$ws.Range("I1").Value = "reward"
$ws.Range("J1").Value = "note"
$ws.Range("K1").Value = "Check"

# ---------------------------------------------------------------------------
# 3. Data rows 2-10.
# ---------------------------------------------------------------------------
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 4
$ws.Range("F2").Value = 3
$ws.Range("G2").Value = -100
$ws.Range("H2").Value = 2
$ws.Range("I2").Value = 100
$ws.Range("J2").Value = "defeat"

$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 500
$ws.Range("H3").Value = 2
$ws.Range("I3").Value = 100
$ws.Range("J3").Value = "victory+3oak+coinx3"

$ws.Range("D4").Value = 3
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 2
$ws.Range("G4").Value = 1300
$ws.Range("H4").Value = 2
$ws.Range("I4").Value = 100
$ws.Range("J4").Value = "victory+3oakx2+coinx5+FC"

$ws.Range("C5").Value = 3
$ws.Range("D5").Value = 2
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = -200
$ws.Range("H5").Value = 3
$ws.Range("I5").Value = 200
$ws.Range("J5").Value = "defeat"

$ws.Range("C6").Value = 2
$ws.Range("D6").Value = 3
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 400
$ws.Range("H6").Value = 3
$ws.Range("I6").Value = 200
$ws.Range("J6").Value = "victory+3oakx2"

$ws.Range("D7").Value = 4
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 2
$ws.Range("G7").Value = 800
$ws.Range("H7").Value = 3
$ws.Range("I7").Value = 200
$ws.Range("J7").Value = "victory+4oak+coinx4"

$ws.Range("C8").Value = 2
$ws.Range("D8").Value = 3
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = -500
$ws.Range("H8").Value = 4
$ws.Range("I8").Value = 500
$ws.Range("J8").Value = "defeat"

$ws.Range("A9").Value = 2
$ws.Range("C9").Value = 2
$ws.Range("D9").Value = 4
$ws.Range("G9").Value = 700
$ws.Range("H9").Value = 4
$ws.Range("I9").Value = 500
$ws.Range("J9").Value = "victory+4oak"

$ws.Range("D10").Value = 5
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 2
$ws.Range("G10").Value = 1800
$ws.Range("H10").Value = 4
$ws.Range("I10").Value = 500
$ws.Range("J10").Value = "victory+5oak+coinx3+FC"

# ---------------------------------------------------------------------------
# 4. Column K "Check" formulas -- K2 stands alone, K3:K10 form a shared
#    formula group (mirrors how the original author filled the column down).
# ---------------------------------------------------------------------------
$ws.Range("K2").Formula = "=SUM(A2:F2)"
$ws.Range("K3:K10").Formula = "=SUM(A3:F3)"

# ---------------------------------------------------------------------------
# 5. Column widths for G:J (values chosen so the engine's internal 1/7-unit
#    rounding lands as close as possible to the original bestFit widths of
#    11.83203125 / 11.83203125 / 14.1640625 / 24.33203125).
# ---------------------------------------------------------------------------
$ws.Columns.Item(7).ColumnWidth = 11.142857142857142
$ws.Columns.Item(8).ColumnWidth = 11.142857142857142
$ws.Columns.Item(9).ColumnWidth = 13.428571428571429
$ws.Columns.Item(10).ColumnWidth = 23.571428571428573

# ---------------------------------------------------------------------------
# 6. Selections on the pre-existing sheets, matching the target workbook.
# ---------------------------------------------------------------------------
$wsNormal = $wb.Worksheets.Item("Normal")
$wsNormal.Range("A1:XFD1").Select() | Out-Null

$wsMonkey = $wb.Worksheets.Item("MonkeyBusiness")
$wsMonkey.Range("H2:H10").Select() | Out-Null

# ---------------------------------------------------------------------------
# 7. Leave SeaBattle the active / selected sheet with H10 as active cell.
# ---------------------------------------------------------------------------
$ws.Select() | Out-Null
$ws.Range("H10").Select() | Out-Null
